$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.100305199623108
$ws.Range("B1").Value = 1.617159605026245
$ws.Range("C1").Value = 2.848811864852905
$ws.Range("D1").Value = 1.513748526573181
$ws.Range("E1").Value = 0.8201173543930054
